$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and a row-order swap for rows 50-51)
$updates = @(
    @{ Cell = 'D2'; Value = '57.906.08' }
    @{ Cell = 'E2'; Value = '  +2.73%  ' }
    @{ Cell = 'D3'; Value = '3.063.34' }
    @{ Cell = 'E3'; Value = '  +2.32%  ' }
    @{ Cell = 'E4'; Value = '  +0.07%  ' }
    @{ Cell = 'D5'; Value = '518.62' }
    @{ Cell = 'E5'; Value = '  +2.47%  ' }
    @{ Cell = 'D6'; Value = '142.35' }
    @{ Cell = 'E6'; Value = '  +3.15%  ' }
    @{ Cell = 'E7'; Value = '  +0.06%  ' }
    @{ Cell = 'D8'; Value = '0.435' }
    @{ Cell = 'E8'; Value = '  +1.25%  ' }
    @{ Cell = 'E9'; Value = '  +1.92%  ' }
    @{ Cell = 'D10'; Value = '0.107' }
    @{ Cell = 'E10'; Value = '  +0.09%  ' }
    @{ Cell = 'E11'; Value = '  +2.91%  ' }
    @{ Cell = 'D12'; Value = '3.590.06' }
    @{ Cell = 'E12'; Value = '  +2.47%  ' }
    @{ Cell = 'E13'; Value = '  +3.34%  ' }
    @{ Cell = 'D14'; Value = '25.80' }
    @{ Cell = 'E14'; Value = '  +0.39%  ' }
    @{ Cell = 'D15'; Value = '0.0000163' }
    @{ Cell = 'E15'; Value = '  +0.16%  ' }
    @{ Cell = 'D16'; Value = '57.931.08' }
    @{ Cell = 'E16'; Value = '  +2.83%  ' }
    @{ Cell = 'D17'; Value = '3.060.40' }
    @{ Cell = 'E17'; Value = '  +2.21%  ' }
    @{ Cell = 'D18'; Value = '6.08' }
    @{ Cell = 'E18'; Value = '  +1.56%  ' }
    @{ Cell = 'D19'; Value = '12.87' }
    @{ Cell = 'E19'; Value = '  -0.57%  ' }
    @{ Cell = 'D20'; Value = '8.14' }
    @{ Cell = 'E20'; Value = '  +0.99%  ' }
    @{ Cell = 'D21'; Value = '331.13' }
    @{ Cell = 'E21'; Value = '  -0.21%  ' }
    @{ Cell = 'D22'; Value = '0.999' }
    @{ Cell = 'E22'; Value = '  -0.18%  ' }
    @{ Cell = 'D23'; Value = '0.499' }
    @{ Cell = 'E23'; Value = '  +0.99%  ' }
    @{ Cell = 'D24'; Value = '65.85' }
    @{ Cell = 'E24'; Value = '  +1.80%  ' }
    @{ Cell = 'D26'; Value = '0.998' }
    @{ Cell = 'E26'; Value = '  -0.13%  ' }
    @{ Cell = 'D27'; Value = '0.0₃0901' }
    @{ Cell = 'E27'; Value = '  -4.90%  ' }
    @{ Cell = 'D28'; Value = '6.41' }
    @{ Cell = 'E28'; Value = '  +0.63%  ' }
    @{ Cell = 'D29'; Value = '7.23' }
    @{ Cell = 'E29'; Value = '  +4.35%  ' }
    @{ Cell = 'D30'; Value = '1.82' }
    @{ Cell = 'E30'; Value = '  +2.47%  ' }
    @{ Cell = 'D31'; Value = '1.19' }
    @{ Cell = 'E31'; Value = '  +3.11%  ' }
    @{ Cell = 'D32'; Value = '20.67' }
    @{ Cell = 'E32'; Value = '  +1.97%  ' }
    @{ Cell = 'D33'; Value = '154.77' }
    @{ Cell = 'E33'; Value = '  +1.39%  ' }
    @{ Cell = 'D34'; Value = '4.52' }
    @{ Cell = 'E34'; Value = '  +0.63%  ' }
    @{ Cell = 'E35'; Value = '  +4.17%  ' }
    @{ Cell = 'D36'; Value = '5.95' }
    @{ Cell = 'E36'; Value = '  +2.57%  ' }
    @{ Cell = 'E37'; Value = '  +0.76%  ' }
    @{ Cell = 'D38'; Value = '0.0675' }
    @{ Cell = 'E38'; Value = '  +2.30%  ' }
    @{ Cell = 'D39'; Value = '3.107.30' }
    @{ Cell = 'E39'; Value = '  +2.54%  ' }
    @{ Cell = 'D40'; Value = '3.92' }
    @{ Cell = 'E40'; Value = '  +3.53%  ' }
    @{ Cell = 'D41'; Value = '36.62' }
    @{ Cell = 'E41'; Value = '  -0.75%  ' }
    @{ Cell = 'E42'; Value = '  +0.09%  ' }
    @{ Cell = 'D43'; Value = '0.655' }
    @{ Cell = 'E43'; Value = '  +0.49%  ' }
    @{ Cell = 'D44'; Value = '2.274.30' }
    @{ Cell = 'E44'; Value = '  +4.03%  ' }
    @{ Cell = 'D45'; Value = '0.0258' }
    @{ Cell = 'E45'; Value = '  +10.15%  ' }
    @{ Cell = 'D46'; Value = '20.76' }
    @{ Cell = 'E46'; Value = '  +6.46%  ' }
    @{ Cell = 'D47'; Value = '1.36' }
    @{ Cell = 'E47'; Value = '  +1.17%  ' }
    @{ Cell = 'D48'; Value = '5.88' }
    @{ Cell = 'E48'; Value = '  +0.72%  ' }
    @{ Cell = 'D49'; Value = '0.930' }
    @{ Cell = 'E49'; Value = '  +0.89%  ' }
    @{ Cell = 'B50'; Value = 'Bittensor' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = 'D50'; Value = '260.67' }
    @{ Cell = 'E50'; Value = '  +12.72%  ' }
    @{ Cell = 'B51'; Value = 'SuiNetwork' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui' }
    @{ Cell = 'D51'; Value = '0.728' }
    @{ Cell = 'E51'; Value = '  +8.21%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = '@'
    $cell.Value = $u.Value
    $cell.Style = 'Normal'
}
